$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "Req Traceability Matrix" cell value changes ---
$ws1.Range("K2").Value = 2
$ws1.Range("L2").Value = 3
$ws1.Range("M2").Value = 5
$ws1.Range("N2").Value = 2

$ws1.Range("B5").Value = 1
$ws1.Range("K5").Value = ""
$ws1.Range("L5").Value = ""

$ws1.Range("B8").Value = 4
$ws1.Range("M8").Value = "x"

$ws1.Range("B9").Value = 7
$ws1.Range("H9").Value = "x"
$ws1.Range("L9").Value = "x"
$ws1.Range("N9").Value = "x"

$ws1.Range("B10").Value = 6
$ws1.Range("H10").Value = "x"
$ws1.Range("J10").Value = "x"
$ws1.Range("L10").Value = "x"
